# Auto-generated edit script: updates market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several sheets,
# reflecting a refreshed price-data snapshot from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2001959.4
$ws.Range("J17").Value = 2501824.2
$ws.Range("L17").Value = 7505472.600000001
$ws.Range("N17").Value = -7505808.600000001
# Row 28
$ws.Range("H28").Value = 1465.25
$ws.Range("J28").Value = 2774.25
$ws.Range("L28").Value = 2774.25
$ws.Range("N28").Value = -3744.25
# Row 40
$ws.Range("H40").Value = 21431964
$ws.Range("J40").Value = 75001870
$ws.Range("L40").Value = 75001870
$ws.Range("N40").Value = -75002220
# Row 70
$ws.Range("H70").Value = 7245.6787
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 9967.315000000001
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 29901.945
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -30441.945
# Row 73
$ws.Range("H73").Value = 7245.6787
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 9967.315000000001
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 29901.945
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -31773.945
# Row 82
$ws.Range("H82").Value = 783.1667
$ws.Range("I82").Value = 783.1667
$ws.Range("K82").Value = 2349.5001
$ws.Range("M82").Value = -1943.5001
# Row 85
$ws.Range("H85").Value = 783.1667
$ws.Range("I85").Value = 783.1667
$ws.Range("K85").Value = 2349.5001
$ws.Range("M85").Value = -945.5001000000002
# Row 101
$ws.Range("H101").Value = 6493873
$ws.Range("I101").Value = 11364028
$ws.Range("K101").Value = 34092084
$ws.Range("M101").Value = -34090462
# Row 107
$ws.Range("H107").Value = 31252446
$ws.Range("I107").Value = 16669275
$ws.Range("K107").Value = 16669275
$ws.Range("M107").Value = -16667355
# Row 131
$ws.Range("H131").Value = 19078.154
$ws.Range("I131").Value = 19084.666
$ws.Range("K131").Value = 57253.99800000001
$ws.Range("M131").Value = -52213.99800000001
# Row 132
$ws.Range("H132").Value = 7597.381
$ws.Range("I132").Value = 2865.75
$ws.Range("J132").Value = 13906.223
$ws.Range("K132").Value = 8597.25
$ws.Range("L132").Value = 41718.669
$ws.Range("M132").Value = -6067.25
$ws.Range("N132").Value = -46778.669
# Row 137
$ws.Range("H137").Value = 26201.844
$ws.Range("J137").Value = 6442.52
$ws.Range("L137").Value = 19327.56
$ws.Range("N137").Value = -24427.56
# Row 138
$ws.Range("H138").Value = 3152.8616
$ws.Range("I138").Value = 1878.3077
$ws.Range("J138").Value = 3471.5
$ws.Range("K138").Value = 5634.9231
$ws.Range("L138").Value = 10414.5
$ws.Range("M138").Value = -494.9231
$ws.Range("N138").Value = -20694.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4765.405
$ws.Range("I61").Value = 1913.3103
$ws.Range("J61").Value = 11127.77
$ws.Range("K61").Value = 1913.3103
$ws.Range("L61").Value = 11127.77
$ws.Range("M61").Value = -1701.3103
$ws.Range("N61").Value = -11551.77
# Row 63
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -2314
# Row 66
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -11568
# Row 134
$ws.Range("H134").Value = 89999
$ws.Range("J134").Value = 89999
$ws.Range("L134").Value = 89999
$ws.Range("N134").Value = -100139
# Row 136
$ws.Range("H136").Value = 4765.405
$ws.Range("I136").Value = 1913.3103
$ws.Range("J136").Value = 11127.77
$ws.Range("K136").Value = 5739.9309
$ws.Range("L136").Value = 33383.31
$ws.Range("M136").Value = -3189.9309
$ws.Range("N136").Value = -38483.31
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20336
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("N53").Value = 0
# Row 134
$ws.Range("H134").Value = 3223.8157
$ws.Range("I134").Value = 2226.2666
$ws.Range("J134").Value = 6964.625
$ws.Range("K134").Value = 6678.7998
$ws.Range("L134").Value = 20893.875
$ws.Range("M134").Value = -4143.7998
$ws.Range("N134").Value = -25963.875

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 464
$ws.Range("I5").Value = 512.5
$ws.Range("J5").Value = 270
$ws.Range("K5").Value = 512.5
$ws.Range("L5").Value = 270
$ws.Range("M5").Value = -400.5
$ws.Range("N5").Value = -494
# Row 26
$ws.Range("H26").Value = 99994
$ws.Range("J26").Value = 99994
$ws.Range("L26").Value = 99994
$ws.Range("N26").Value = -100568
# Row 31
$ws.Range("H31").Value = 186432.42
$ws.Range("I31").Value = 228224.27
$ws.Range("J31").Value = 2548.3
$ws.Range("K31").Value = 228224.27
$ws.Range("L31").Value = 2548.3
$ws.Range("M31").Value = -227929.27
$ws.Range("N31").Value = -3138.3
# Row 32
$ws.Range("H32").Value = 17123.75
$ws.Range("J32").Value = 19997.5
$ws.Range("L32").Value = 19997.5
$ws.Range("N32").Value = -20629.5
# Row 34
$ws.Range("H34").Value = 186432.42
$ws.Range("I34").Value = 228224.27
$ws.Range("J34").Value = 2548.3
$ws.Range("K34").Value = 228224.27
$ws.Range("L34").Value = 2548.3
$ws.Range("M34").Value = -228022.27
$ws.Range("N34").Value = -2952.3
# Row 62
$ws.Range("H62").Value = 7194.3335
$ws.Range("J62").Value = 7133.2
$ws.Range("L62").Value = 7133.2
$ws.Range("N62").Value = -8381.200000000001
# Row 65
$ws.Range("H65").Value = 7194.3335
$ws.Range("J65").Value = 7133.2
$ws.Range("L65").Value = 35666
$ws.Range("N65").Value = -41906
# Row 105
$ws.Range("H105").Value = 5190.25
$ws.Range("I105").Value = 1973.125
$ws.Range("J105").Value = 7763.95
$ws.Range("K105").Value = 1973.125
$ws.Range("L105").Value = 7763.95
$ws.Range("M105").Value = -226.125
$ws.Range("N105").Value = -11257.95
# Row 122
$ws.Range("H122").Value = 3247.2222
$ws.Range("I122").Value = 1882.4
$ws.Range("J122").Value = 4953.25
$ws.Range("K122").Value = 5647.200000000001
$ws.Range("L122").Value = 14859.75
$ws.Range("M122").Value = -3197.200000000001
$ws.Range("N122").Value = -19759.75
# Row 132
$ws.Range("H132").Value = 4808.95
$ws.Range("I132").Value = 3083.1035
$ws.Range("J132").Value = 9358.909
$ws.Range("K132").Value = 9249.3105
$ws.Range("L132").Value = 28076.727
$ws.Range("M132").Value = -6719.3105
$ws.Range("N132").Value = -33136.727
# Row 134
$ws.Range("H134").Value = 4383.2544
$ws.Range("I134").Value = 5007.609
$ws.Range("J134").Value = 2174
$ws.Range("K134").Value = 15022.827
$ws.Range("L134").Value = 6522
$ws.Range("M134").Value = -12487.827
$ws.Range("N134").Value = -11592

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 122154.836
$ws.Range("I4").Value = 28714.363
$ws.Range("K4").Value = 86143.08900000001
$ws.Range("M4").Value = -86031.08900000001
# Row 98
$ws.Range("H98").Value = 543.7692
$ws.Range("I98").Value = 465.66666
$ws.Range("J98").Value = 610.7143
$ws.Range("K98").Value = 1396.99998
$ws.Range("L98").Value = 1832.1429
$ws.Range("M98").Value = 101.0000199999999
$ws.Range("N98").Value = -4828.1429
# Row 131
$ws.Range("H131").Value = 6850899.5
$ws.Range("I131").Value = 38462550
$ws.Range("J131").Value = 1708.5
$ws.Range("K131").Value = 115387650
$ws.Range("L131").Value = 5125.5
$ws.Range("M131").Value = -115382610
$ws.Range("N131").Value = -15205.5
# Row 136
$ws.Range("H136").Value = 5188.75
$ws.Range("I136").Value = 5188.75
$ws.Range("K136").Value = 15566.25
$ws.Range("M136").Value = -10466.25

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 11522.714
$ws.Range("J126").Value = 4165
$ws.Range("L126").Value = 12495
$ws.Range("N126").Value = -17435
# Row 132
$ws.Range("H132").Value = 45925.043
$ws.Range("I132").Value = 56930.156
$ws.Range("J132").Value = 4105.6
$ws.Range("K132").Value = 170790.468
$ws.Range("L132").Value = 12316.8
$ws.Range("M132").Value = -168260.468
$ws.Range("N132").Value = -17376.8

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 7611.8
$ws.Range("I132").Value = 7284.6313
$ws.Range("J132").Value = 8647.833000000001
$ws.Range("K132").Value = 21853.8939
$ws.Range("L132").Value = 25943.499
$ws.Range("M132").Value = -19323.8939
$ws.Range("N132").Value = -31003.499

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 20840948
$ws.Range("I62").Value = 7598
$ws.Range("K62").Value = 7598
$ws.Range("M62").Value = -6974
# Row 65
$ws.Range("H65").Value = 20840948
$ws.Range("I65").Value = 7598
$ws.Range("K65").Value = 37990
$ws.Range("M65").Value = -34870
# Row 68
$ws.Range("H68").Value = 19623
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 19623
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 136
$ws.Range("H136").Value = 257326.72
$ws.Range("I136").Value = 286631.53
$ws.Range("K136").Value = 859894.5900000001
$ws.Range("M136").Value = -857344.5900000001
